$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4325566589832306
$ws.Range("B1").Value = 1.875313997268677
$ws.Range("C1").Value = 2.204512596130371
$ws.Range("D1").Value = 1.938581228256226
$ws.Range("E1").Value = 0.9797055721282959
